# time_log.xlsx edit: "Finish when will i retire"
#
# Adds a 5th table column ("Milestones"), splits the combined
# "Finish exercises + .5 hours of cheatsheets / JS100 = 37 hours" note
# into a Notes cell + a Milestones cell, moves the "LS95 = 14 hours"
# milestone note into the new Milestones column, and records one more
# small problem finished on the last day (Hours 1.75 -> 2, note updated).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Add the new "Milestones" table column -----------------------------
$tbl = $ws.ListObjects.Item(1)
$null = $tbl.ListColumns.Add()
$ws.Range("E1").Value = "Milestones"

# Data cells default to center alignment (matches the table's dxf for this
# column), the totals row cell is left aligned like the rest of the row.
$ws.Range("E2:E64").HorizontalAlignment = -4108   # xlCenter
$ws.Range("E65").HorizontalAlignment = -4131      # xlLeft

# --- 2. Move the "LS95 = 14 hours" milestone from Notes to Milestones -----
$ws.Range("D9").Value = ""
$ws.Range("E9").Value = "LS95 = 14 hours"
$ws.Range("E9").HorizontalAlignment = -4131       # xlLeft, matches old D9 style

# --- 3. Split the combined note on 44479 into Notes + Milestones ----------
$ws.Range("D41").Value = "Finish exercises + .5 hours of cheatsheets"
$ws.Range("E41").Value = "JS100 = 37 hours"

# --- 4. Update the last logged day: one more small problem finished -------
$ws.Range("C60").Value = 2
$ws.Range("D60").Value = "Finished 4 small problems"

# --- 5. Cosmetics: column widths + view position ---------------------------
$ws.Range("D1").EntireColumn.ColumnWidth = 55.28515625 - 5/6
$ws.Range("E1").EntireColumn.ColumnWidth = 15.42578125 - 5/6

$excel.ActiveWindow.ScrollRow = 39
$excel.ActiveWindow.ScrollColumn = 1
$null = $ws.Range("D43").Select()
